# Update the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (commit: "Updated cryptos list on Wed May 29 21:51:30 UTC 2024
# with GitHub Actions").
#
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
# Most rows only need their Price (D) and Volume(1h) (E) values refreshed.
# A couple of rows (24/25 and 39/40) were re-ranked, so the Coin name, link
# and price/volume moved to a different row than before.
#
# All of these columns hold plain text in the source workbook (prices like
# "67.732.11" or "0.0000276" are strings, not numbers). Excel will happily
# auto-convert a numeric-looking string into a real number when you assign
# it through .Value, which would corrupt values such as "0.0000276"
# (rendered as "2.76E-05") or strip formatting. To avoid that we briefly
# force the cell's number format to Text ("@") before assigning the value,
# then restore the cell style to Normal so no stray formatting is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$CellRef, [string]$Text)
    $c = $Sheet.Range($CellRef)
    $c.NumberFormat = "@"
    $c.Value = $Text
    $c.Style = "Normal"
}

# Rows whose Coin (B) and Link (C) swapped along with their Price/Volume.
Set-TextValue $ws "B24" "PEPE"
Set-TextValue $ws "C24" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws "D24" "0.0000149"
Set-TextValue $ws "E24" "  -7.71%  "

Set-TextValue $ws "B25" "Litecoin"
Set-TextValue $ws "C25" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws "D25" "83.63"
Set-TextValue $ws "E25" "  +0.35%  "

Set-TextValue $ws "B39" "Mantle"
Set-TextValue $ws "C39" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D39" "1.01"
Set-TextValue $ws "E39" "  -1.58%  "

Set-TextValue $ws "B40" "Kaspa"
Set-TextValue $ws "C40" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D40" "0.138"
Set-TextValue $ws "E40" "  -1.42%  "

# Remaining rows: refresh Price (D) and/or Volume(1h) (E) only.
$updates = @(
    @{ Row = 2;  D = "67.732.11";  E = "  -0.99%  " },
    @{ Row = 3;  D = "3.786.35";   E = "  -1.45%  " },
    @{ Row = 4;  D = "0.997";      E = "  -0.18%  " },
    @{ Row = 5;  D = "597.12";     E = "  -0.91%  " },
    @{ Row = 6;  D = "169.72";     E = "  +0.31%  " },
    @{ Row = 7;  D = "3.784.87";   E = "  -1.42%  " },
    @{ Row = 8;  E = "  -0.10%  " },
    @{ Row = 9;  D = "0.524";      E = "  -0.67%  " },
    @{ Row = 10; D = "0.165";      E = "  -0.95%  " },
    @{ Row = 11; E = "  -0.37%  " },
    @{ Row = 12; E = "  -1.11%  " },
    @{ Row = 13; D = "0.0000276";  E = "  +2.54%  " },
    @{ Row = 14; D = "36.53";      E = "  -1.58%  " },
    @{ Row = 15; D = "4.408.13";   E = "  -1.69%  " },
    @{ Row = 16; D = "3.775.76";   E = "  -1.93%  " },
    @{ Row = 17; D = "18.66";      E = "  +0.89%  " },
    @{ Row = 18; D = "67.624.30";  E = "  -1.20%  " },
    @{ Row = 19; D = "7.19";       E = "  -2.83%  " },
    @{ Row = 20; E = "  +0.83%  " },
    @{ Row = 21; D = "10.58";      E = "  -4.29%  " },
    @{ Row = 22; D = "468.31";     E = "  -0.40%  " },
    @{ Row = 23; D = "0.720";      E = "  -1.92%  " },
    @{ Row = 26; D = "2.21";       E = "  -0.87%  " },
    @{ Row = 27; D = "12.15";      E = "  +0.03%  " },
    @{ Row = 28; D = "10.31";      E = "  +2.36%  " },
    @{ Row = 29; E = "  +0.07%  " },
    @{ Row = 30; E = "  -2.02%  " },
    @{ Row = 31; D = "3.921.94";   E = "  -1.70%  " },
    @{ Row = 32; D = "7.64";       E = "  -0.59%  " },
    @{ Row = 33; D = "30.55";      E = "  -3.25%  " },
    @{ Row = 34; D = "2.23";       E = "  -3.64%  " },
    @{ Row = 35; D = "9.12";       E = "  -2.72%  " },
    @{ Row = 36; D = "3.736.26";   E = "  -1.81%  " },
    @{ Row = 37; D = "3.83";       E = "  +3.36%  " },
    @{ Row = 38; D = "0.104";      E = "  -1.11%  " },
    @{ Row = 41; D = "5.81";       E = "  -2.19%  " },
    @{ Row = 42; E = "  -0.06%  " },
    @{ Row = 43; D = "0.312";      E = "  -0.93%  " },
    @{ Row = 45; D = "8.70";       E = "  -0.22%  " },
    @{ Row = 46; D = "1.94";       E = "  -2.06%  " },
    @{ Row = 47; D = "45.82";      E = "  -2.75%  " },
    @{ Row = 48; D = "396.95";     E = "  -4.32%  " },
    @{ Row = 49; D = "0.000271";   E = "  -7.02%  " },
    @{ Row = 50; D = "139.78";     E = "  -1.25%  " },
    @{ Row = 51; D = "0.0353";     E = "  -2.14%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        Set-TextValue $ws "D$($u.Row)" $u.D
    }
    Set-TextValue $ws "E$($u.Row)" $u.E
}
